$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("E1").Value = "failed?"
$ws.Range("F1").Value = "description"
$ws.Range("G1").Value = "TODO after"

# --- Mark every existing data row (2-10) as "removed" in column E ---
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 5).Value = "removed"
}

# --- Row 8 & 9 gain a description / TODO after note ---
$ws.Range("F8").Value = "Renewed version"
$ws.Range("G8").Value = "annotation, DEG, pseudotime"
$ws.Range("F9").Value = "Renewed version"
$ws.Range("G9").Value = "annotation, DEG, pseudotime"

# --- Row 10 previously had no Name (col B); now it does ---
$ws.Range("B10").Value = "2022-06-10 16-29-41"

# --- New rows 11-13: individual corrected pipeline runs ---
$ws.Range("A11").Value = "results"
$ws.Range("B11").Value = "Pipe_SCTv2_corrected_13-06"
$ws.Range("C11").Value = "individual"
$ws.Range("D11").Value = "SCTv2 corrected BL_C"
$ws.Range("F11").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G11").Value = "integration (old/new selection), annotation, DEG, pseudotime"

$ws.Range("A12").Value = "results"
$ws.Range("B12").Value = "Pipe_SCTv2_corrected_13-06"
$ws.Range("C12").Value = "individual"
$ws.Range("D12").Value = "SCTv2 corrected BL_A"
$ws.Range("F12").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G12").Value = "integration (old/new selection), annotation, DEG, pseudotime"

$ws.Range("A13").Value = "results"
$ws.Range("B13").Value = "Pipe_SCTv2_corrected_13-06"
$ws.Range("C13").Value = "individual"
$ws.Range("D13").Value = "SCTv2 corrected BL_N"
$ws.Range("F13").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G13").Value = "integration (old/new selection), annotation, DEG, pseudotime"

# --- New rows 14-17: integration results for new/old selections ---
$ws.Range("A14").Value = "results"
$ws.Range("B14").Value = "2022-06-13 13-32-07"
$ws.Range("C14").Value = "integration"
$ws.Range("D14").Value = "SCTv2 corrected BL_A + BL_C new selection"
$ws.Range("F14").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G14").Value = "annotation, DEG, pseudotime"

$ws.Range("A15").Value = "results"
$ws.Range("B15").Value = "2022-06-13 13-33-22"
$ws.Range("C15").Value = "integration"
$ws.Range("D15").Value = "SCTv2 corrected BL_A + BL_C old selection"
$ws.Range("F15").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G15").Value = "annotation, DEG, pseudotime"

$ws.Range("A16").Value = "results"
$ws.Range("B16").Value = "2022-06-13 13-34-02"
$ws.Range("C16").Value = "integration"
$ws.Range("D16").Value = "SCTv2 corrected BL_N + BL_C new selection"
$ws.Range("F16").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G16").Value = "annotation, DEG, pseudotime"

$ws.Range("A17").Value = "results"
$ws.Range("B17").Value = "2022-06-13 13-35-10"
$ws.Range("C17").Value = "integration"
$ws.Range("D17").Value = "SCTv2 corrected BL_N + BL_C old selection"
$ws.Range("F17").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G17").Value = "annotation, DEG, pseudotime"

# --- Column widths for the two new columns (best-fit, matching F/G content) ---
$ws.Columns.Item(6).ColumnWidth = 67.45182291666667
$ws.Columns.Item(7).ColumnWidth = 56.592447916666664

# --- Selection / active cell matches the post-edit state ---
$ws.Range("B20").Select()
